$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.061.28'
$ws.Range("E2").Value = '  +1.13%  '

$ws.Range("D3").Value = '1.849.97'
$ws.Range("E3").Value = '  +2.30%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = '''235.44'
$ws.Range("E5").Value = '  +2.12%  '

$ws.Range("E6").Value = '  +2.75%  '

$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("D8").Value = '''42.44'
$ws.Range("E8").Value = '  +8.61%  '

$ws.Range("E9").Value = '  +2.00%  '

$ws.Range("D10").Value = '''0.0694'
$ws.Range("E10").Value = '  +2.15%  '

$ws.Range("D11").Value = '''0.0987'
$ws.Range("E11").Value = '  -0.69%  '

$ws.Range("D12").Value = '2.115.49'
$ws.Range("E12").Value = '  +2.04%  '

$ws.Range("D13").Value = '1.861.84'
$ws.Range("E13").Value = '  +2.85%  '

$ws.Range("D14").Value = '''11.38'
$ws.Range("E14").Value = '  +1.82%  '

$ws.Range("D15").Value = '''0.677'
$ws.Range("E15").Value = '  +2.16%  '

$ws.Range("E16").Value = '  +2.66%  '

$ws.Range("D17").Value = '35.015.46'
$ws.Range("E17").Value = '  +0.97%  '

$ws.Range("D18").Value = '''70.07'
$ws.Range("E18").Value = '  +1.10%  '

$ws.Range("D19").Value = '0.0₃0795'
$ws.Range("E19").Value = '  +1.66%  '

$ws.Range("D20").Value = '''241.18'
$ws.Range("E20").Value = '  +0.62%  '

$ws.Range("D21").Value = '''12.16'
$ws.Range("E21").Value = '  +2.61%  '

$ws.Range("D22").Value = '''4.79'
$ws.Range("E22").Value = '  +2.87%  '

$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").Value = '''2.26'
$ws.Range("E24").Value = '  +1.10%  '

$ws.Range("D25").Value = '''171.30'
$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("D26").Value = '''1.85'
$ws.Range("E26").Value = '  +23.19%  '

$ws.Range("D27").Value = '''7.90'
$ws.Range("E27").Value = '  +1.93%  '

$ws.Range("D28").Value = '''17.66'
$ws.Range("E28").Value = '  +2.85%  '

$ws.Range("E29").Value = '  +3.36%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '''0.0558'
$ws.Range("E30").Value = '  +2.42%  '

$ws.Range("B31").Value = 'BinanceUSD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D31").Value = '''1.01'
$ws.Range("E31").Value = '  -0.26%  '

$ws.Range("D32").Value = '''3.99'
$ws.Range("E32").Value = '  -2.01%  '

$ws.Range("D33").Value = '''3.97'
$ws.Range("E33").Value = '  +1.24%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = '''2.03'
$ws.Range("E34").Value = '  +14.11%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '''1.64'
$ws.Range("E35").Value = '  +23.41%  '

$ws.Range("D36").Value = '''0.772'
$ws.Range("E36").Value = '  +11.05%  '

$ws.Range("E37").Value = '  -1.75%  '

$ws.Range("E38").Value = '  +11.85%  '

$ws.Range("D39").Value = '''91.73'
$ws.Range("E39").Value = '  +0.29%  '

$ws.Range("E40").Value = '  +5.09%  '

$ws.Range("D41").Value = '1.348.45'
$ws.Range("E41").Value = '  +1.63%  '

$ws.Range("D42").Value = '''15.04'
$ws.Range("E42").Value = '  +5.41%  '

$ws.Range("D43").Value = '''12.94'
$ws.Range("E43").Value = '  +84.99%  '

$ws.Range("E44").Value = '  +6.28%  '

$ws.Range("E45").Value = '  -3.38%  '

$ws.Range("E46").Value = '  +2.43%  '

$ws.Range("D47").Value = '''6.41'
$ws.Range("E47").Value = '  +3.08%  '

$ws.Range("E48").Value = '  +3.82%  '

$ws.Range("D49").Value = '2.026.59'
$ws.Range("E49").Value = '  +1.36%  '

$ws.Range("D50").Value = '''3.47'
$ws.Range("E50").Value = '  +17.28%  '

$ws.Range("E51").Value = '  +1.82%  '
